# Bondtech Kit reference added
# Rows 3-8 (Vendor = Bondtech) of the Dragon hotend/extruder BOM previously
# pointed at six different individual Bondtech product pages. They now all
# point at the single "BMG Internals Set for HextrudORT" kit page, the
# Part Description is unified, Make/Buy becomes "(BUY) KIT", and the
# bearing row's QTY note becomes "2*" (it used to be "1*").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://www.bondtech.se/product/bmg-internals-set-for-hextrudort/"
$newDesc = "Included in BMG Internals Set for HextrudORT"
$newMakeBuy = "(BUY) KIT"

$rows = 3,4,5,6,7,8

foreach ($r in $rows) {
    $ws.Range("F$r").Value = $newDesc
    $ws.Range("G$r").Value = $newMakeBuy
}

# Bearing row: QTY note goes from "1*" to "2*"
$ws.Range("H8").Value = "2*"

# Re-point all six Vendor URL hyperlinks (K3:K8) at the new kit page,
# replacing the previous individual-product links.
$ws.Hyperlinks.Delete()

foreach ($r in $rows) {
    $ws.Hyperlinks.Add($ws.Range("K$r"), $newUrl, "", "", $newUrl)
}

# Restore the selection left behind by the edit (diff shows activeCell moved
# from O7 to E3).
$ws.Range("E3").Select()
